$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer (index 1 in Word's Footers collection) holds the Pearson logo
# whose docPr id="2" currently reads "image1.png" -> should become "image2.png"
$f1 = $sec.Footers.Item(1)
$shp1 = $f1.Range.InlineShapes.Item(1)
$shp1.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# Footer (index 2) holds the Pearson logo whose docPr id="3" currently
# reads "image1.png" -> should become "image2.png"
$f2 = $sec.Footers.Item(2)
$shp2 = $f2.Range.InlineShapes.Item(1)
$shp2.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# Header (index 2) holds the BTEC logo whose docPr id="1" currently
# reads "image2.jpg" -> should become "image1.jpg"
$h1 = $sec.Headers.Item(2)
$shp3 = $h1.Range.InlineShapes.Item(1)
$shp3.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
